$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.170.79"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").Value = "1.846.70"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7020"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.52"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3045"
$ws.Range("E8").Value = "  -3.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07406"
$ws.Range("E9").Value = "  +3.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.37"
$ws.Range("E10").Value = "  -5.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08129"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.867.20"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7265"
$ws.Range("E13").Value = "  -3.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.211"
$ws.Range("E14").Value = "  -3.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "88.76"
$ws.Range("E15").Value = "  -4.20%  "
$ws.Range("D16").Value = "29.145.09"
$ws.Range("E16").Value = "  -2.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "238.13"
$ws.Range("E18").Value = "  -4.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.06"
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "2.091.24"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.584"
$ws.Range("E24").Value = "  -3.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.996"
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.59"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1451"
$ws.Range("E27").Value = "  -7.39%  "
$ws.Range("E28").Value = "  -3.23%  "
$ws.Range("E29").Value = "  -3.69%  "
$ws.Range("E30").Value = "  -5.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.516"
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.490"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.993"
$ws.Range("E33").Value = "  -4.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05178"
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.184"
$ws.Range("E35").Value = "  -5.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.032"
$ws.Range("E36").Value = "  +3.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7028"
$ws.Range("E37").Value = "  -8.52%  "
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01868"
$ws.Range("E39").Value = "  -4.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.673"
$ws.Range("E40").Value = "  -3.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9490"
$ws.Range("E41").Value = "  +7.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.008"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("E43").Value = "  -6.10%  "
$ws.Range("D44").Value = "1.068.35"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.88"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.741"
$ws.Range("E48").Value = "  -6.01%  "
$ws.Range("D49").Value = "1.988.20"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.023"
$ws.Range("E50").Value = "  -6.99%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.124"
$ws.Range("E51").Value = "  -4.33%  "
